# Update the three worksheets (Lightweight, Featherweight, Heavyweight) so
# that the "Age" column becomes a "Win via" column with fight-result data.

$wb = $excel.ActiveWorkbook

# --- Lightweight sheet ---
$ws = $wb.Worksheets.Item("Lightweight")
$ws.Range("C1").Value = "Win via"

$ws.Range("A2").Value = "Dagestan"
$ws.Range("B2").Value = "Khabib"
$ws.Range("C2").Value = "SUB"

$ws.Range("A3").Value = "USA"
$ws.Range("B3").Value = "Justin"
$ws.Range("C3").Value = "KO"

$ws.Range("A4").Value = "Brazil"
$ws.Range("B4").Value = "Charles"
$ws.Range("C4").Value = "TKO"

$ws.Range("C13").Select()

# --- Featherweight sheet ---
$ws = $wb.Worksheets.Item("Featherweight")
$ws.Range("C1").Value = "Win via"

$ws.Range("A2").Value = "Ireland"
$ws.Range("B2").Value = "Connor"
$ws.Range("C2").Value = "KO"

$ws.Range("A3").Value = "USA"
$ws.Range("B3").Value = "Tony"
$ws.Range("C3").Value = "SUB"

$ws.Range("A4").Value = "Louisiana"
$ws.Range("B4").Value = "Dustin"
$ws.Range("C4").Value = "TKO"

$ws.Range("F16").Select()

# --- Heavyweight sheet ---
$ws = $wb.Worksheets.Item("Heavyweight")
$ws.Range("C1").Value = "Win via"

$ws.Range("A2").Value = "South Africa"
$ws.Range("B2").Value = "Francis"
$ws.Range("C2").Value = "KO"

$ws.Range("A3").Value = "France"
$ws.Range("B3").Value = "Gane"
$ws.Range("C3").Value = "SUB"

$ws.Range("A4").Value = "USA"
$ws.Range("B4").Value = "Daniel"
$ws.Range("C4").Value = "KO"

$ws.Range("C4").Select()
